$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 368.57144
$ws.Range("I12").Value = 92.5
$ws.Range("J12").Value = 479
$ws.Range("K12").Value = 92.5
$ws.Range("L12").Value = 479
$ws.Range("M12").Value = 77.5
$ws.Range("N12").Value = -819

$ws.Range("H33").Value = 3860.182
$ws.Range("I33").Value = 4134
$ws.Range("K33").Value = 4134
$ws.Range("M33").Value = -3905

$ws.Range("H100").Value = 5428.2144
$ws.Range("I100").Value = 1713
$ws.Range("J100").Value = 10381.833
$ws.Range("K100").Value = 1713
$ws.Range("L100").Value = 10381.833
$ws.Range("M100").Value = -1172
$ws.Range("N100").Value = -11463.833

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 44124836
$ws.Range("I61").Value = 71434280
$ws.Range("J61").Value = 25008230
$ws.Range("K61").Value = 71434280
$ws.Range("L61").Value = 25008230
$ws.Range("M61").Value = -71434068
$ws.Range("N61").Value = -25008654

$ws.Range("H74").Value = 17334484
$ws.Range("I74").Value = 22728388
$ws.Range("K74").Value = 22728388
$ws.Range("M74").Value = -22727514

$ws.Range("H77").Value = 17334484
$ws.Range("I77").Value = 22728388
$ws.Range("K77").Value = 113641940
$ws.Range("M77").Value = -113637572

$ws.Range("H102").Value = 12470.083
$ws.Range("I102").Value = 10460.111
$ws.Range("K102").Value = 10460.111
$ws.Range("M102").Value = -8838.111000000001

$ws.Range("H109").Value = 42190.332
$ws.Range("J109").Value = 42190.332
$ws.Range("L109").Value = 42190.332
$ws.Range("N109").Value = -44964.332

$ws.Range("H114").Value = 57494.75
$ws.Range("J114").Value = 57494.75
$ws.Range("L114").Value = 57494.75
$ws.Range("N114").Value = -66172.75

$ws.Range("H119").Value = 63938
$ws.Range("J119").Value = 63938
$ws.Range("L119").Value = 63938
$ws.Range("N119").Value = -73614

$ws.Range("H122").Value = 3191.5833
$ws.Range("I122").Value = 2106.8667
$ws.Range("J122").Value = 4999.4443
$ws.Range("K122").Value = 6320.6001
$ws.Range("L122").Value = 14998.3329
$ws.Range("M122").Value = -3870.6001
$ws.Range("N122").Value = -19898.3329

$ws.Range("H136").Value = 44124836
$ws.Range("I136").Value = 71434280
$ws.Range("J136").Value = 25008230
$ws.Range("K136").Value = 214302840
$ws.Range("L136").Value = 75024690
$ws.Range("M136").Value = -214300290
$ws.Range("N136").Value = -75029790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 59999.688
$ws.Range("J115").Value = 59999.688
$ws.Range("L115").Value = 59999.688
$ws.Range("N115").Value = -63133.688

$ws.Range("H134").Value = 6253984
$ws.Range("I134").Value = 3703.818
$ws.Range("K134").Value = 11111.454
$ws.Range("M134").Value = -8576.454000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1229984.1
$ws.Range("I31").Value = 4131
$ws.Range("J31").Value = 1374202.1
$ws.Range("K31").Value = 4131
$ws.Range("L31").Value = 1374202.1
$ws.Range("M31").Value = -3836
$ws.Range("N31").Value = -1374792.1

$ws.Range("H34").Value = 1229984.1
$ws.Range("I34").Value = 4131
$ws.Range("J34").Value = 1374202.1
$ws.Range("K34").Value = 4131
$ws.Range("L34").Value = 1374202.1
$ws.Range("M34").Value = -3929
$ws.Range("N34").Value = -1374606.1

$ws.Range("H86").Value = 5207.952
$ws.Range("I86").Value = 5207.2856
$ws.Range("J86").Value = 5209.2856
$ws.Range("K86").Value = 5207.2856
$ws.Range("L86").Value = 5209.2856
$ws.Range("M86").Value = -4084.2856
$ws.Range("N86").Value = -7455.2856

$ws.Range("H89").Value = 5207.952
$ws.Range("I89").Value = 5207.2856
$ws.Range("J89").Value = 5209.2856
$ws.Range("K89").Value = 26036.428
$ws.Range("L89").Value = 26046.428
$ws.Range("M89").Value = -20420.428
$ws.Range("N89").Value = -37278.428

$ws.Range("H92").Value = 64663.332
$ws.Range("J92").Value = 64663.332
$ws.Range("L92").Value = 64663.332
$ws.Range("N92").Value = -69655.33199999999

$ws.Range("H140").Value = 72050
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2705.3333
$ws.Range("I12").Value = 3261.1428
$ws.Range("K12").Value = 9783.428400000001
$ws.Range("M12").Value = -9610.428400000001

$ws.Range("H23").Value = 917.875
$ws.Range("I23").Value = 1650.3334
$ws.Range("K23").Value = 4951.0002
$ws.Range("M23").Value = -4716.0002

$ws.Range("H33").Value = 152.35294
$ws.Range("I33").Value = 44.444443
$ws.Range("J33").Value = 273.75
$ws.Range("K33").Value = 266.666658
$ws.Range("L33").Value = 1642.5
$ws.Range("M33").Value = 16.33334200000002
$ws.Range("N33").Value = -2208.5

$ws.Range("H59").Value = 1600
$ws.Range("I59").Value = 1600
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 4800
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -4260
$ws.Range("N59").ClearContents()

$ws.Range("H114").Value = 1561
$ws.Range("J114").Value = 999.6667
$ws.Range("L114").Value = 2999.0001
$ws.Range("N114").Value = -9507.000100000001

$ws.Range("H134").Value = 9105.096
$ws.Range("I134").Value = 1371.75
$ws.Range("J134").Value = 12198.434
$ws.Range("K134").Value = 4115.25
$ws.Range("L134").Value = 36595.302
$ws.Range("M134").Value = 954.75
$ws.Range("N134").Value = -46735.302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 15628.5
$ws.Range("I13").Value = 16378.75
$ws.Range("J13").Value = 14128
$ws.Range("K13").Value = 16378.75
$ws.Range("L13").Value = 14128
$ws.Range("M13").Value = -16239.75
$ws.Range("N13").Value = -14406

$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H41").Value = 3875
$ws.Range("I41").Value = 2666.6667
$ws.Range("J41").Value = 7500
$ws.Range("K41").Value = 2666.6667
$ws.Range("L41").Value = 7500
$ws.Range("M41").Value = -2311.6667
$ws.Range("N41").Value = -8210

$ws.Range("H70").Value = 12420.923
$ws.Range("I70").Value = 17496.285
$ws.Range("J70").Value = 6499.6665
$ws.Range("K70").Value = 17496.285
$ws.Range("L70").Value = 6499.6665
$ws.Range("M70").Value = -17226.285
$ws.Range("N70").Value = -7039.6665

$ws.Range("H73").Value = 12420.923
$ws.Range("I73").Value = 17496.285
$ws.Range("J73").Value = 6499.6665
$ws.Range("K73").Value = 17496.285
$ws.Range("L73").Value = 6499.6665
$ws.Range("M73").Value = -16560.285
$ws.Range("N73").Value = -8371.666499999999

$ws.Range("H104").Value = 76045
$ws.Range("J104").Value = 76045
$ws.Range("L104").Value = 76045
$ws.Range("N104").Value = -83033

$ws.Range("H111").Value = 86194.75
$ws.Range("J111").Value = 86194.75
$ws.Range("L111").Value = 86194.75
$ws.Range("N111").Value = -92328.75

$ws.Range("H113").Value = 5002.6
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 5003.25
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 5003.25
$ws.Range("M113").Value = -2830
$ws.Range("N113").Value = -9343.25

$ws.Range("H121").Value = 39995
$ws.Range("J121").Value = 39995
$ws.Range("L121").Value = 39995
$ws.Range("N121").Value = -43489

$ws.Range("H122").Value = 2200.0715
$ws.Range("I122").Value = 1822.7778
$ws.Range("J122").Value = 2879.2
$ws.Range("K122").Value = 5468.3334
$ws.Range("L122").Value = 8637.599999999999
$ws.Range("M122").Value = -3018.3334
$ws.Range("N122").Value = -13537.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4576.48
$ws.Range("I40").Value = 3815.2144
$ws.Range("J40").Value = 5545.364
$ws.Range("K40").Value = 3815.2144
$ws.Range("L40").Value = 5545.364
$ws.Range("M40").Value = -3679.2144
$ws.Range("N40").Value = -5817.364

$ws.Range("H46").Value = 3242.077
$ws.Range("J46").Value = 4976.5
$ws.Range("L46").Value = 4976.5
$ws.Range("N46").Value = -5352.5

$ws.Range("H61").Value = 1579.4
$ws.Range("I61").Value = 1579.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1579.4
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1377.4
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 1579.4
$ws.Range("I113").Value = 1579.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1579.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 590.5999999999999
$ws.Range("N113").ClearContents()

$ws.Range("H136").Value = 77785.19
$ws.Range("I136").Value = 12589.546
$ws.Range("K136").Value = 37768.638
$ws.Range("M136").Value = -35218.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 747.2222
$ws.Range("I113").Value = 597
$ws.Range("J113").Value = 1137.8
$ws.Range("K113").Value = 1791
$ws.Range("L113").Value = 3413.4
$ws.Range("M113").Value = 379
$ws.Range("N113").Value = -7753.4

$ws.Range("H122").Value = 2214.2163
$ws.Range("I122").Value = 2265.3547
$ws.Range("K122").Value = 6796.0641
$ws.Range("M122").Value = -4346.0641
